$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on the Price/Volume columns being updated so that
# numeric-looking strings (e.g. "1.001", "316.60") are stored as literal text,
# matching the source data which keeps these as inline strings, not numbers.
$textCells = @(
    "D2", "E2", "D3", "E3", "E4", "D5", "E5", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "D21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D27", "E27", "B28", "C28", "D28", "E28", "B29", "C29", "D29", "E29", "D30", "E30", "D31", "E31", "D32", "E32", "D33", "E33", "E34", "D35", "E35", "D36", "E36", "D37", "E37", "D38", "E38", "E39", "D40", "E40", "E41", "D42", "E42", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "D48", "E48", "D49", "E49", "D50", "D51", "E51"
)
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range("D2").Value = "28.613.08"
$ws.Range("E2").Value = "  +1.46%  "
$ws.Range("D3").Value = "1.829.59"
$ws.Range("E3").Value = "  +1.34%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "316.60"
$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").Value = "0.5337"
$ws.Range("E7").Value = "  -0.87%  "
$ws.Range("D8").Value = "0.3998"
$ws.Range("E8").Value = "  +5.58%  "
$ws.Range("D9").Value = "0.07776"
$ws.Range("E9").Value = "  +3.83%  "
$ws.Range("D10").Value = "1.122"
$ws.Range("E10").Value = "  +2.07%  "
$ws.Range("E11").Value = "  +0.10%  "
$ws.Range("D12").Value = "21.24"
$ws.Range("E12").Value = "  +3.21%  "
$ws.Range("D13").Value = "6.334"
$ws.Range("E13").Value = "  +1.95%  "
$ws.Range("D14").Value = "7.605"
$ws.Range("E14").Value = "  +2.83%  "
$ws.Range("E15").Value = "  +0.14%  "
$ws.Range("D16").Value = "1.828.08"
$ws.Range("E16").Value = "  +1.36%  "
$ws.Range("D17").Value = "93.25"
$ws.Range("E17").Value = "  +3.60%  "
$ws.Range("D18").Value = "0.00001095"
$ws.Range("E18").Value = "  +2.64%  "
$ws.Range("D19").Value = "0.06594"
$ws.Range("E19").Value = "  +1.34%  "
$ws.Range("D20").Value = "17.84"
$ws.Range("E20").Value = "  +2.28%  "
$ws.Range("D21").Value = "1.001"
$ws.Range("D22").Value = "6.107"
$ws.Range("E22").Value = "  +2.82%  "
$ws.Range("D23").Value = "28.626.53"
$ws.Range("E23").Value = "  +1.42%  "
$ws.Range("D24").Value = "11.22"
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").Value = "2.240"
$ws.Range("E25").Value = "  +7.25%  "
$ws.Range("D26").Value = "20.86"
$ws.Range("E26").Value = "  +1.47%  "
$ws.Range("D27").Value = "157.07"
$ws.Range("E27").Value = "  +0.64%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "2.439"
$ws.Range("E28").Value = "  +4.37%  "
$ws.Range("B29").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C29").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D29").Value = "2.037.19"
$ws.Range("E29").Value = "  +1.28%  "
$ws.Range("D30").Value = "125.72"
$ws.Range("E30").Value = "  +2.86%  "
$ws.Range("D31").Value = "1.166"
$ws.Range("E31").Value = "  +3.38%  "
$ws.Range("D32").Value = "0.1130"
$ws.Range("E32").Value = "  +1.18%  "
$ws.Range("D33").Value = "5.767"
$ws.Range("E33").Value = "  +2.72%  "
$ws.Range("E34").Value = "  -0.28%  "
$ws.Range("D35").Value = "0.07395"
$ws.Range("E35").Value = "  +6.01%  "
$ws.Range("D36").Value = "0.2279"
$ws.Range("E36").Value = "  +2.04%  "
$ws.Range("D37").Value = "0.02359"
$ws.Range("E37").Value = "  +2.33%  "
$ws.Range("D38").Value = "8.951"
$ws.Range("E38").Value = "  +5.62%  "
$ws.Range("E39").Value = "  +2.34%  "
$ws.Range("D40").Value = "11.44"
$ws.Range("E40").Value = "  +2.15%  "
$ws.Range("E41").Value = "  +2.04%  "
$ws.Range("D42").Value = "1.199"
$ws.Range("E42").Value = "  +1.71%  "
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("D44").Value = "1.393"
$ws.Range("E44").Value = "  -2.53%  "
$ws.Range("D45").Value = "13.57"
$ws.Range("E45").Value = "  +1.00%  "
$ws.Range("D46").Value = "0.5957"
$ws.Range("E46").Value = "  +3.00%  "
$ws.Range("D47").Value = "3.716"
$ws.Range("E47").Value = "  +0.76%  "
$ws.Range("D48").Value = "125.87"
$ws.Range("E48").Value = "  +0.33%  "
$ws.Range("D49").Value = "2.005"
$ws.Range("E49").Value = "  +3.69%  "
$ws.Range("D50").Value = "1.195"
$ws.Range("D51").Value = "0.06971"
$ws.Range("E51").Value = "  +2.10%  "

# Restore the default cell style (the text number-format above is a transient
# aid for entry; the resulting cells should keep the workbook's original style).
foreach ($ref in $textCells) {
    $ws.Range($ref).Style = "Normal"
}
